# "ADD ALL ARTICLE AT A TIME CHANE THE CLICK  METHOD"
# The Addarticle sheet held a bunch of placeholder/test "article" names in
# A2:A18 (dfs, dfsdf, dfsdgg, ... , "TV 50\"") left over from manual testing of
# the add-article form. They get retyped in one pass with a fresh batch of
# placeholder values, and the selection left on A13 (scrolled so row 10 is at
# the top) instead of A17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Addarticle")
$ws.Activate()

$newArticles = @(
    "FASFA",
    "ASFG",
    "GASG",
    "ASGS",
    "AGE",
    "GSDE",
    "TE",
    "G",
    "SDFSDFS",
    "DFSDFSD",
    "SDFSDFSSDF",
    "SDF",
    "BXCBXCB",
    "CB",
    "XCBXC",
    "BXCBXCB",
    "XCBXC"
)

for ($i = 0; $i -lt $newArticles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newArticles[$i]
}

# Move the view the way the user left it: scrolled down with A13 selected.
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
